$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1380.9375
$ws.Cells.Item(17, 10).Value = 1315.6
$ws.Cells.Item(17, 12).Value = 3946.8
$ws.Cells.Item(17, 14).Value = -4282.799999999999
$ws.Cells.Item(51, 8).Value = 9666.666999999999
$ws.Cells.Item(51, 9).Value = 9500
$ws.Cells.Item(51, 11).Value = 9500
$ws.Cells.Item(51, 13).Value = -9016
$ws.Cells.Item(62, 8).Value = 2833.1667
$ws.Cells.Item(62, 9).Value = 2399.8
$ws.Cells.Item(62, 11).Value = 2399.8
$ws.Cells.Item(62, 13).Value = -1775.8
$ws.Cells.Item(65, 8).Value = 2833.1667
$ws.Cells.Item(65, 9).Value = 2399.8
$ws.Cells.Item(65, 11).Value = 11999
$ws.Cells.Item(65, 13).Value = -8879
$ws.Cells.Item(74, 8).Value = 10777.333
$ws.Cells.Item(74, 9).Value = 9399.200000000001
$ws.Cells.Item(74, 10).Value = 12500
$ws.Cells.Item(74, 11).Value = 9399.200000000001
$ws.Cells.Item(74, 12).Value = 12500
$ws.Cells.Item(74, 13).Value = -8463.200000000001
$ws.Cells.Item(74, 14).Value = -14372
$ws.Cells.Item(77, 8).Value = 10777.333
$ws.Cells.Item(77, 9).Value = 9399.200000000001
$ws.Cells.Item(77, 10).Value = 12500
$ws.Cells.Item(77, 11).Value = 46996
$ws.Cells.Item(77, 12).Value = 62500
$ws.Cells.Item(77, 13).Value = -42316
$ws.Cells.Item(77, 14).Value = -71860
$ws.Cells.Item(137, 8).Value = 3141.7727
$ws.Cells.Item(137, 9).Value = 2970.95
$ws.Cells.Item(137, 11).Value = 8912.849999999999
$ws.Cells.Item(137, 13).Value = -6362.849999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 35740.777
$ws.Cells.Item(2, 9).Value = 39708.375
$ws.Cells.Item(2, 11).Value = 39708.375
$ws.Cells.Item(2, 13).Value = -39595.375
$ws.Cells.Item(32, 8).Value = 16133349
$ws.Cells.Item(32, 9).Value = 16671064
$ws.Cells.Item(32, 11).Value = 16671064
$ws.Cells.Item(32, 13).Value = -16670777
$ws.Cells.Item(44, 8).Value = 37999.5
$ws.Cells.Item(44, 10).Value = 37999.5
$ws.Cells.Item(44, 12).Value = 37999.5
$ws.Cells.Item(44, 14).Value = -38975.5
$ws.Cells.Item(45, 8).Value = 2992.5715
$ws.Cells.Item(45, 9).Value = 2589.8
$ws.Cells.Item(45, 11).Value = 2589.8
$ws.Cells.Item(45, 13).Value = -2212.8
$ws.Cells.Item(55, 8).Value = 450000
$ws.Cells.Item(55, 10).Value = 450000
$ws.Cells.Item(55, 12).Value = 450000
$ws.Cells.Item(55, 14).Value = -450630
$ws.Cells.Item(61, 8).Value = 1372.6842
$ws.Cells.Item(61, 9).Value = 1282.8334
$ws.Cells.Item(61, 11).Value = 1282.8334
$ws.Cells.Item(61, 13).Value = -1070.8334
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 35740.777
$ws.Cells.Item(116, 9).Value = 39708.375
$ws.Cells.Item(116, 11).Value = 39708.375
$ws.Cells.Item(116, 13).Value = -37414.375
$ws.Cells.Item(119, 8).Value = 86666
$ws.Cells.Item(119, 10).Value = 86666
$ws.Cells.Item(119, 12).Value = 86666
$ws.Cells.Item(119, 14).Value = -96342
$ws.Cells.Item(122, 8).Value = 2780
$ws.Cells.Item(122, 9).Value = 2502.5
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 7507.5
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -5057.5
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(132, 8).Value = 2062.5557
$ws.Cells.Item(132, 9).Value = 2104.9033
$ws.Cells.Item(132, 11).Value = 6314.7099
$ws.Cells.Item(132, 13).Value = -3784.7099
$ws.Cells.Item(136, 8).Value = 1372.6842
$ws.Cells.Item(136, 9).Value = 1282.8334
$ws.Cells.Item(136, 11).Value = 3848.5002
$ws.Cells.Item(136, 13).Value = -1298.5002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 35740.777
$ws.Cells.Item(3, 9).Value = 39708.375
$ws.Cells.Item(3, 11).Value = 39708.375
$ws.Cells.Item(3, 13).Value = -39594.375
$ws.Cells.Item(33, 8).Value = 2000
$ws.Cells.Item(33, 9).Value = 2000
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 13).Value = -1664
$ws.Cells.Item(86, 8).Value = 2556.75
$ws.Cells.Item(86, 9).Value = 2340.5334
$ws.Cells.Item(86, 11).Value = 2340.5334
$ws.Cells.Item(86, 13).Value = -1217.5334
$ws.Cells.Item(89, 8).Value = 2556.75
$ws.Cells.Item(89, 9).Value = 2340.5334
$ws.Cells.Item(89, 11).Value = 11702.667
$ws.Cells.Item(89, 13).Value = -6086.666999999999
$ws.Cells.Item(134, 8).Value = 1149.0256
$ws.Cells.Item(134, 9).Value = 953.1111
$ws.Cells.Item(134, 11).Value = 2859.3333
$ws.Cells.Item(134, 13).Value = -324.3332999999998
$ws.Cells.Item(137, 8).Value = 69143.63
$ws.Cells.Item(137, 10).Value = 69143.63
$ws.Cells.Item(137, 12).Value = 69143.63
$ws.Cells.Item(137, 14).Value = -79343.63

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1805.9286
$ws.Cells.Item(31, 9).Value = 1798.7407
$ws.Cells.Item(31, 10).Value = 2000
$ws.Cells.Item(31, 11).Value = 1798.7407
$ws.Cells.Item(31, 12).Value = 2000
$ws.Cells.Item(31, 13).Value = -1503.7407
$ws.Cells.Item(31, 14).Value = -2590
$ws.Cells.Item(34, 8).Value = 1805.9286
$ws.Cells.Item(34, 9).Value = 1798.7407
$ws.Cells.Item(34, 10).Value = 2000
$ws.Cells.Item(34, 11).Value = 1798.7407
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = -1596.7407
$ws.Cells.Item(34, 14).Value = -2404
$ws.Cells.Item(86, 8).Value = 19085.633
$ws.Cells.Item(86, 10).Value = 10749.6
$ws.Cells.Item(86, 12).Value = 10749.6
$ws.Cells.Item(86, 14).Value = -12995.6
$ws.Cells.Item(89, 8).Value = 19085.633
$ws.Cells.Item(89, 10).Value = 10749.6
$ws.Cells.Item(89, 12).Value = 53748
$ws.Cells.Item(89, 14).Value = -64980
$ws.Cells.Item(132, 8).Value = 2896.3076
$ws.Cells.Item(132, 9).Value = 2896.3076
$ws.Cells.Item(132, 11).Value = 8688.9228
$ws.Cells.Item(132, 13).Value = -6158.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 904.6667
$ws.Cells.Item(9, 10).Value = 857
$ws.Cells.Item(9, 12).Value = 2571
$ws.Cells.Item(9, 14).Value = -3019
$ws.Cells.Item(88, 8).Value = 4214.143
$ws.Cells.Item(88, 9).Value = 1500
$ws.Cells.Item(88, 10).Value = 4666.5
$ws.Cells.Item(88, 11).Value = 4500
$ws.Cells.Item(88, 12).Value = 13999.5
$ws.Cells.Item(88, 13).Value = -4072
$ws.Cells.Item(88, 14).Value = -14855.5
$ws.Cells.Item(91, 8).Value = 4214.143
$ws.Cells.Item(91, 9).Value = 1500
$ws.Cells.Item(91, 10).Value = 4666.5
$ws.Cells.Item(91, 11).Value = 4500
$ws.Cells.Item(91, 12).Value = 13999.5
$ws.Cells.Item(91, 13).Value = -3018
$ws.Cells.Item(91, 14).Value = -16963.5
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents()
$ws.Cells.Item(137, 8).Value = 2553.1304
$ws.Cells.Item(137, 9).Value = 2214.8
$ws.Cells.Item(137, 10).Value = 3187.5
$ws.Cells.Item(137, 11).Value = 6644.400000000001
$ws.Cells.Item(137, 12).Value = 9562.5
$ws.Cells.Item(137, 13).Value = -1544.400000000001
$ws.Cells.Item(137, 14).Value = -19762.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5501.5557
$ws.Cells.Item(126, 9).Value = 5185
$ws.Cells.Item(126, 10).Value = 6134.6665
$ws.Cells.Item(126, 11).Value = 15555
$ws.Cells.Item(126, 12).Value = 18403.9995
$ws.Cells.Item(126, 13).Value = -13085
$ws.Cells.Item(126, 14).Value = -23343.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 1534.75
$ws.Cells.Item(32, 9).Value = 1534.75
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 1534.75
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -1217.75
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 4022.111
$ws.Cells.Item(61, 9).Value = 2243
$ws.Cells.Item(61, 11).Value = 2243
$ws.Cells.Item(61, 13).Value = -2041
$ws.Cells.Item(82, 8).Value = 1677.2858
$ws.Cells.Item(82, 9).Value = 1613.6666
$ws.Cells.Item(82, 11).Value = 1613.6666
$ws.Cells.Item(82, 13).Value = -1252.6666
$ws.Cells.Item(85, 8).Value = 1677.2858
$ws.Cells.Item(85, 9).Value = 1613.6666
$ws.Cells.Item(85, 11).Value = 1613.6666
$ws.Cells.Item(85, 13).Value = -365.6666
$ws.Cells.Item(113, 8).Value = 4022.111
$ws.Cells.Item(113, 9).Value = 2243
$ws.Cells.Item(113, 11).Value = 2243
$ws.Cells.Item(113, 13).Value = -73

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1228.8928
$ws.Cells.Item(126, 9).Value = 1110.7142
$ws.Cells.Item(126, 10).Value = 1583.4286
$ws.Cells.Item(126, 11).Value = 3332.1426
$ws.Cells.Item(126, 12).Value = 4750.2858
$ws.Cells.Item(126, 13).Value = -862.1425999999997
$ws.Cells.Item(126, 14).Value = -9690.2858
